# Hortaliza, Feria Lagunitas de Puerto Montt - Ají
# A new weekly price observation was inserted as row 292 (pushing the
# existing rows 292:378 down to 293:379). The dimension grows from
# A1:R378 to A1:R379.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 292, shifting everything below
# it (including formatting) down by one row.
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row 292 with the new observation.
$ws.Range("A292").Value = 4
$ws.Range("B292").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C292").Value = "Los Lagos"
$ws.Range("D292").Value = 44985
$ws.Range("E292").Value = 10
$ws.Range("F292").Value = 100112021
$ws.Range("G292").Value = "Ají"
$ws.Range("H292").Value = "Inferno"
$ws.Range("I292").Value = "Primera"
$ws.Range("J292").Value = 180
$ws.Range("K292").Value = 18000
$ws.Range("L292").Value = 20000
$ws.Range("M292").Value = 19000
$ws.Range("N292").Value = "$/caja 10 kilos"
$ws.Range("O292").Value = "Región de Arica y Parinacota"
$ws.Range("P292").Value = 1900
$ws.Range("Q292").Value = 10
$ws.Range("R292").Value = "Hortaliza"
